$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-20 Friday" "2024-12-21 Saturday"

Replace-Text "59×15=" "59×93="
Replace-Text "38×55=" "28×57="
Replace-Text "35×20=" "60×15="
Replace-Text "57×81=" "28×38="
Replace-Text "58×99=" "18×35="
Replace-Text "19×23=" "60×91="
Replace-Text "50×80=" "14×53="
Replace-Text "58×37=" "90×18="
Replace-Text "18×34=" "46×14="
Replace-Text "31×52=" "32×76="
Replace-Text "97×62=" "64×29="
Replace-Text "64×51=" "50×83="
Replace-Text "97×92=" "15×96="
Replace-Text "63×47=" "21×47="
Replace-Text "94×49=" "89×91="
Replace-Text "94×67=" "95×22="
Replace-Text "96×12=" "87×12="
Replace-Text "73×76=" "72×11="
Replace-Text "43×80=" "68×95="
Replace-Text "49×16=" "35×24="
Replace-Text "93×15=" "97×70="
Replace-Text "58×94=" "74×44="
Replace-Text "94×60=" "82×76="
Replace-Text "33×57=" "17×18="
Replace-Text "90×14=" "97×21="
